$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Eliptik eğri" row (row 11: "Eliptik eğri" / "245,0603999895975").
# Deleting the entire row shifts the rows below (DSA, IDEA) up by one.
$ws.Rows.Item(11).Delete()

# Reflect the cell selection left behind in the saved file.
$ws.Range("B18").Select()
